$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 2 is the handed-back file (228b4934-...md)
$wsZh.Range("B2").Value = "Handed back: in sync with en-US"
$wsZh.Range("D2").Value = "2016-02-22 13:59:58"
$wsZh.Range("E2").Value = "228b4934-faed-44a0-b362-1b99a5cea0b9.md"
$wsZh.Range("F2").Value = "228b4934-faed-44a0-b362-1b99a5cea0b9.c7f65fee7b20e509e2ecb2f5a389c22b40dbfd31.zh-cn.xlf"
$wsZh.Range("G2").Value = "2016-02-22 14:00:58"

# de-de sheet: row 2 is the handed-back file (228b4934-...md)
$wsDe.Range("B2").Value = "Handed back: in sync with en-US"
$wsDe.Range("D2").Value = "2016-02-22 14:00:15"
$wsDe.Range("E2").Value = "228b4934-faed-44a0-b362-1b99a5cea0b9.md"
$wsDe.Range("F2").Value = "228b4934-faed-44a0-b362-1b99a5cea0b9.c7f65fee7b20e509e2ecb2f5a389c22b40dbfd31.de-de.xlf"
$wsDe.Range("G2").Value = "2016-02-22 14:01:41"
